# Edit: Tue, Jun 23, 2020  7:04:46 AM
#
# 1) Re-point the three data tables (on slides 14, 15 and 16) from the
#    custom "Table_0" table style onto the built-in
#    "{A009BB01-7BA1-4339-A30B-D27909667DC9}" table style.
# 2) Swap the presentation's colour theme from the "Integral" (Red
#    Violet) palette over to the stock Office palette.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Table styles
# ---------------------------------------------------------------
$newTableStyle = "{A009BB01-7BA1-4339-A30B-D27909667DC9}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)

    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle, $true)
        }
    }
}

# ---------------------------------------------------------------
# 2) Theme colours -> Office palette (BGR-packed long values, as
#    used by the RGB COM property)
# ---------------------------------------------------------------
$officeColors = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
